# Add a new "trend_summary" worksheet as the last sheet in the workbook.
# It brings together, side by side, the "trend" column from each of the
# five per-metric trend tables (dph, dps, n_clicks, n_encounters, n_trains).

$wb = $excel.ActiveWorkbook

$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$ws.Name = "trend_summary"

$dph = $wb.Worksheets.Item("trend_table_dph")
$dps = $wb.Worksheets.Item("trend_table_dps")
$nclicks = $wb.Worksheets.Item("trend_table_n_clicks")
$nenc = $wb.Worksheets.Item("trend_table_n_encounters")
$ntrains = $wb.Worksheets.Item("trend_table_n_trains")

# Row/group labels (col A) + station labels (col B) are identical across
# every trend table, so just bring them over from the first one, blank A1
# included.
$dph.Range("A1:B5").Copy($ws.Range("A1:B5"))

# New header labels for the five metric columns.
$ws.Range("C1").Value = "dph"
$ws.Range("D1").Value = "dps"
$ws.Range("E1").Value = "n_clicks"
$ws.Range("F1").Value = "n_encounters"
$ws.Range("G1").Value = "n_trains"

# Pull each metric's "trend" column (column C in the source sheets) into
# its own column here.
$dph.Range("C2:C5").Copy($ws.Range("C2:C5"))
$dps.Range("C2:C5").Copy($ws.Range("D2:D5"))
$nclicks.Range("C2:C5").Copy($ws.Range("E2:E5"))
$nenc.Range("C2:C5").Copy($ws.Range("F2:F5"))
$ntrains.Range("C2:C5").Copy($ws.Range("G2:G5"))
